$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 16
$ws.Range("B2").Value = 88
$ws.Range("C2").Value = 47217
$ws.Range("D2").Value = "layak"
$ws.Range("E2").Value = 0.52
$ws.Range("F2").Value = 75

# Row 3
$ws.Range("A3").Value = 19
$ws.Range("B3").Value = 81
$ws.Range("C3").Value = 46450
$ws.Range("D3").Value = "layak"
$ws.Range("E3").Value = 0.76
$ws.Range("F3").Value = 75

# Row 4
$ws.Range("A4").Value = 87
$ws.Range("B4").Value = 78
$ws.Range("C4").Value = 47995
$ws.Range("D4").Value = "layak"
$ws.Range("E4").Value = 0.8663333333333333
$ws.Range("F4").Value = 75

# Row 5 (A5, B5, C5 unchanged)
$ws.Range("D5").Value = "layak"
$ws.Range("E5").Value = 0.7249333333333333
$ws.Range("F5").Value = 75

# Row 6
$ws.Range("A6").Value = 88
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 35304
$ws.Range("D6").Value = "layak"
$ws.Range("E6").Value = 0.9696
$ws.Range("F6").Value = 75
